$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix number format on H184:I184 so they match the rest of the row
#    (1 decimal place, same as the other monthly figures) instead of the old
#    2-decimal-place format.
# ---------------------------------------------------------------------------
$ws.Range("H184:I184").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 2. Append the latest monthly data (Apr, May, Jun 2025) that was missing,
#    reusing row 184's formatting (date style + number style + borders) for
#    the newly populated rows.
# ---------------------------------------------------------------------------
$ws.Range("A184:J184").Copy()
$ws.Range("A185:J187").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newMonths = @(
    @(45748, 111.3, 16.2, 7.3, 21,   23.9, 7.6, 0.1, 20.8, 3.4),
    @(45778, 114.4, 10.1, 5.4, 20.2, 25.3, 7.6, 0.1, 21.5, 3.5),
    @(45809, 113.5, 5.6,  4.3, 19.8, 33.1, 7.4, 0.1, 21.5, 3.7)
)

$cols = @("A","B","C","D","E","F","G","H","I","J")
$r = 185
foreach ($monthRow in $newMonths) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $monthRow[$i]
    }
    $r++
}

# ---------------------------------------------------------------------------
# 3. Header row (row 1) re-wrapped to a shorter height.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 43.5

# ---------------------------------------------------------------------------
# 4. Update the view: scroll the frozen pane down and move the selection to
#    the new last data row.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 173
$win.ScrollColumn = 2
$ws.Range("M174").Select()
